$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Email" column: header in H1, value in H2
$ws.Range("H1").Value = "Email"
$ws.Range("H2").Value = "adrianrentea01@gmail.com"

# Size the new column to fit its content (matches the author's bestFit width of 24.5)
$ws.Columns.Item(8).ColumnWidth = 23.67

# Scroll the view over so column B becomes the left-most visible column
$excel.ActiveWindow.ScrollColumn = 2

# Mirror the author's final selection covering the newly-added column
$ws.Range("G1:H2").Select()
